$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.93
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.36
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 2.08
$ws.Range("S4").Value = 2.75
$ws.Range("T4").Value = 1.4
$ws.Range("G5").Value = 1.05
$ws.Range("H5").Value = 9.25
$ws.Range("I5").Value = 21
$ws.Range("J5").Value = 1.27
$ws.Range("K5").Value = 3.6
$ws.Range("L5").Value = 14.5
$ws.Range("Q5").Value = 1.19
$ws.Range("R5").Value = 4.3
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.02
$ws.Range("Y5").Value = 13
$ws.Range("Z5").Value = 7.2
$ws.Range("AA5").Value = 14
$ws.Range("AB5").Value = 5.8
$ws.Range("AC5").Value = 11
$ws.Range("AE5").Value = 29
$ws.Range("AF5").Value = 24
$ws.Range("AI5").Value = 80
$ws.Range("AJ5").Value = 250
$ws.Range("AK5").Value = 80
$ws.Range("AM5").Value = 400
$ws.Range("AN5").Value = 250
$ws.Range("K7").Value = 3.75
$ws.Range("N7").Value = 29
$ws.Range("Q7").Value = 1.22
$ws.Range("R7").Value = 4.33
$ws.Range("W7").Value = 1.95
$ws.Range("X7").Value = 1.8
$ws.Range("Y7").Value = 15
$ws.Range("AO7").Value = 301
$ws.Range("I8").Value = 2.7
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.33
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 2.63
$ws.Range("T8").Value = 1.44
$ws.Range("U8").Value = 1.33
$ws.Range("V8").Value = 3.25
$ws.Range("W8").Value = 1.53
$ws.Range("X8").Value = 2.38
$ws.Range("Y8").Value = 11
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 19
$ws.Range("AL8").Value = 29
$ws.Range("AO8").Value = 126
$ws.Range("G9").Value = 1.33
$ws.Range("I9").Value = 7.5
$ws.Range("J9").Value = 1.8
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 19
$ws.Range("Y9").Value = 8.5
$ws.Range("AA9").Value = 9
$ws.Range("AB9").Value = 9
$ws.Range("AI9").Value = 21
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 81
$ws.Range("AN9").Value = 41
$ws.Range("AO9").Value = 251

Write-Host "Applied odds updates"
